$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11: category, recommendation, new_safety_stock, new_reorder_point,
# new_optimal_inventory, new_holding_cost, potential_saving
$data = @(
    ,@(2, 'garden_tools', 'Giảm Safety Stock từ 1552 → 1241 và Reorder Point từ 189523 → 170570 để tiết kiệm chi phí.', 1241, 170570, 171811, 2144201280, -1190354880)
    ,@(3, 'watches_gifts', 'Giảm Safety Stock từ 501 → 400 và Reorder Point từ 40167 → 36150 để tiết kiệm chi phí.', 400, 36150, 36550, 427635000, -237308760)
    ,@(4, 'furniture_decor', 'Giảm Safety Stock từ 550 → 440 và Reorder Point từ 8075 → 7267 để tiết kiệm chi phí.', 440, 7267, 7707, 78148980, -43165980)
    ,@(5, 'bed_bath_table', 'Giảm Safety Stock từ 146 → 116 và Reorder Point từ 7536 → 6782 để tiết kiệm chi phí.', 116, 6782, 6898, 71739200, -39782080)
    ,@(6, 'electronics', 'Giảm Safety Stock từ 152 → 121 và Reorder Point từ 6341 → 5706 để tiết kiệm chi phí.', 121, 5706, 5827, 59085780, -32750172)
    ,@(7, 'auto', 'Giảm Safety Stock từ 70 → 56 và Reorder Point từ 5684 → 5115 để tiết kiệm chi phí.', 56, 5115, 5171, 44367180, -24619452)
    ,@(8, 'sports_leisure', 'Giảm Safety Stock từ 894 → 715 và Reorder Point từ 3548 → 3193 để tiết kiệm chi phí.', 715, 3193, 3908, 38611040, -21056256)
    ,@(9, 'baby', 'Giảm Safety Stock từ 0 → 0 và Reorder Point từ 3749 → 3374 để tiết kiệm chi phí.', 0, 3374, 3374, 32457880, -18031728)
    ,@(10, 'health_beauty', 'Giảm Safety Stock từ 283 → 226 và Reorder Point từ 2317 → 2085 để tiết kiệm chi phí.', 226, 2085, 2311, 22231820, -12227020)
    ,@(11, 'construction_tools_construction', 'Giảm Safety Stock từ 0 → 0 và Reorder Point từ 2035 → 1831 để tiết kiệm chi phí.', 0, 1831, 1831, 15233920, -8461440)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
